$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell holds a plain-text (inline string) value. Excel
# auto-detects numeric-looking text on assignment and would silently
# convert it to a real number, so we force the cell to Text format first,
# write the literal string, then reset the style back to Normal so no
# stray formatting/style is introduced.
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "68.353.67"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = "  +1.03%  "
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "3.748.72"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = "  -0.82%  "
$cell.Style = "Normal"
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = "  -0.05%  "
$cell.Style = "Normal"
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "595.37"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = "  -0.24%  "
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "167.08"
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = "  -1.29%  "
$cell.Style = "Normal"
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "3.746.05"
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = "  -0.87%  "
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = "  -0.08%  "
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.520"
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = "  -0.88%  "
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.160"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = "  -3.23%  "
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = "  -0.22%  "
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = "  -0.91%  "
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000259"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = "  -7.42%  "
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "36.04"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = "  -1.26%  "
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.375.99"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = "  -0.92%  "
$cell.Style = "Normal"
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "3.747.97"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = "  -1.06%  "
$cell.Style = "Normal"
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "68.350.19"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = "  +1.01%  "
$cell.Style = "Normal"
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "17.91"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = "  -3.32%  "
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "7.01"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = "  -2.31%  "
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = "  -0.23%  "
$cell.Style = "Normal"
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "10.72"
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = "  +1.85%  "
$cell.Style = "Normal"
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "467.62"
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = "  -0.23%  "
$cell.Style = "Normal"
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "0.698"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = "  -2.89%  "
$cell.Style = "Normal"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "84.45"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = "  +0.74%  "
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = "  -3.14%  "
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = "  -0.61%  "
$cell.Style = "Normal"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "11.99"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = "  -1.42%  "
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "10.11"
$cell.Style = "Normal"
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = "  -1.80%  "
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = "  -0.06%  "
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "3.894.79"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = "  -0.94%  "
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = "  -4.64%  "
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = "  -4.21%  "
$cell.Style = "Normal"
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "29.85"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = "  -2.23%  "
$cell.Style = "Normal"
$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = "  -2.27%  "
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = "  +0.61%  "
$cell.Style = "Normal"
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "3.703.81"
$cell.Style = "Normal"
$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = "  -1.08%  "
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = "  -2.67%  "
$cell.Style = "Normal"
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "3.39"
$cell.Style = "Normal"
$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = "  -10.04%  "
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = "  +0.83%  "
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = "  -0.41%  "
$cell.Style = "Normal"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "5.80"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = "  -0.01%  "
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = "  -0.06%  "
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = "  -2.30%  "
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = "  -1.10%  "
$cell.Style = "Normal"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "43.09"
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = "  +9.77%  "
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = "  -1.26%  "
$cell.Style = "Normal"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "45.91"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = "  +0.09%  "
$cell.Style = "Normal"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "146.66"
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = "  +4.19%  "
$cell.Style = "Normal"
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "389.55"
$cell.Style = "Normal"
$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = "  -1.64%  "
$cell.Style = "Normal"
